# Update gh-pages to output generated at 456a3b4
# Applies numeric "想去人数" (want-to-go count) bumps and two cover-image URL
# updates across the 展览, 演出, and 全部类型 worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 37766
$ws.Range("F4").Value  = 640
$ws.Range("F5").Value  = 784
$ws.Range("F6").Value  = 486
$ws.Range("F9").Value  = 857
$ws.Range("F10").Value = 97
$ws.Range("F11").Value = 729
$ws.Range("F12").Value = 565
$ws.Range("F13").Value = 66
$ws.Range("F18").Value = 476
$ws.Range("F23").Value = 2555
$ws.Range("F24").Value = 1040
$ws.Range("F29").Value = 799
$ws.Range("I29").Value = "//i2.hdslb.com/bfs/openplatform/202404/QrBvxNAX1712126496119.jpeg"
$ws.Range("F30").Value = 72
$ws.Range("F31").Value = 1171

# ---- Sheet: 演出 -----------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value  = 420
$ws.Range("F5").Value  = 4
$ws.Range("F10").Value = 14
$ws.Range("F12").Value = 11

# ---- Sheet: 本地生活 ---------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 647

# ---- Sheet: 全部类型 ---------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 647
$ws.Range("F3").Value  = 37766
$ws.Range("F5").Value  = 640
$ws.Range("F6").Value  = 784
$ws.Range("F7").Value  = 486
$ws.Range("F11").Value = 420
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = 857
$ws.Range("F16").Value = 97
$ws.Range("F17").Value = 729
$ws.Range("F18").Value = 565
$ws.Range("F19").Value = 66
$ws.Range("F24").Value = 14
$ws.Range("F29").Value = 476
$ws.Range("F34").Value = 2555
$ws.Range("F35").Value = 1040
$ws.Range("F40").Value = 11
$ws.Range("F41").Value = 799
$ws.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202404/QrBvxNAX1712126496119.jpeg"
$ws.Range("F42").Value = 72
$ws.Range("F43").Value = 1171
